$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add two new worksheets after the existing "States" sheet: "Sheet1"
# (the active/selected tab, holding the new fish Type/Variety data) and
# "Sheet2" (still empty - a placeholder for more Gauge test data).
# ---------------------------------------------------------------------
$statesSheet = $wb.Worksheets.Item("States")

$ws1 = $wb.Worksheets.Add($null, $statesSheet)
$ws1.Name = "Sheet1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# Populate Sheet1 with the fish Type / Variety table.
# ---------------------------------------------------------------------
$data = @(
    @("Type", "Variety"),
    @("Trout", "Rainbow"),
    @("Trout", "Cutthroat"),
    @("Trout", "Brook"),
    @("Trout", "Brown"),
    @("Trout", "Sea"),
    @("Trout", "Lake"),
    @("Trout", "Spake"),
    @("Trout", "Bull"),
    @("Trout", "Golden"),
    @("Trout", "Tiger"),
    @("Trout", "Dolly Vardon"),
    @("Bass", "Smallmouth"),
    @("Bass", "Largemouth"),
    @("Bass", "Peacock"),
    @("Bass", "Striped"),
    @("Bass", "White"),
    @("Bass", "Black"),
    @("Bass", "Yellow")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $data[$i][0]
    $ws1.Cells.Item($row, 2).Value = $data[$i][1]
}

# Cursor ended up one row below the last entry, on Sheet1.
$ws1.Range("A20").Select() | Out-Null

# Sheet2 stays empty; leave its selection at the default A1.

# Make Sheet1 the active tab (matches activeTab="1" / tabSelected on Sheet1).
$ws1.Activate()

# ---------------------------------------------------------------------
# States sheet: a handful of rows were manually resized (no longer
# driven purely by autofit), shrinking from either 28.8 -> 16.8 or
# 43.2 -> 28.8 points.
# ---------------------------------------------------------------------
$rowHeights = @{
    24 = 16.8
    25 = 16.8
    30 = 28.8
    31 = 16.8
    32 = 16.8
    37 = 16.8
    39 = 16.8
    40 = 16.8
    42 = 16.8
    43 = 16.8
    48 = 16.8
    49 = 16.8
    50 = 16.8
}

foreach ($r in $rowHeights.Keys) {
    $statesSheet.Rows.Item($r).RowHeight = $rowHeights[$r]
}
